$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the shared date formula in column B down to the new row and
# populate the new row (13) with this week's final timesheet entry.
$ws.Range("B12").AutoFill($ws.Range("B12:B13"), 0)

$ws.Range("A13").Value = 42904
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = "Finished integrating the properties file into the code base. Introduced some errors along the way that took a while to fix. Also changed glitter to appear only in the cell with the gold in keeping with the standard wumpus world model. "

$ws.Rows.Item(13).RowHeight = 42.75
